$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F304").Value = 5688
$ws.Range("F309").Value = 74292
$ws.Range("F314").Value = 63432
$ws.Range("F321").Value = 90525
$ws.Range("G321").Value = 2792
$ws.Range("F322").Value = 107318
$ws.Range("F324").Value = 232606
$ws.Range("F325").Value = 730318
$ws.Range("G325").Value = 6032
$ws.Range("F326").Value = 426005
$ws.Range("G326").Value = 3744
$ws.Range("F330").Value = 70707
$ws.Range("F331").Value = 150108
$ws.Range("F332").Value = 421479
$ws.Range("G332").Value = 4091
$ws.Range("F333").Value = 258458
$ws.Range("G333").Value = 2787
$ws.Range("F334").Value = 202101
$ws.Range("G334").Value = 3379
$ws.Range("F336").Value = 100223
$ws.Range("F337").Value = 101961
$ws.Range("G337").Value = 2882
$ws.Range("F338").Value = 215655
$ws.Range("G338").Value = 3061
$ws.Range("F339").Value = 600001
$ws.Range("G339").Value = 5193
$ws.Range("F340").Value = 335192
$ws.Range("G340").Value = 3029
$ws.Range("F341").Value = 394112
$ws.Range("G341").Value = 4605
$ws.Range("F342").Value = 185871
$ws.Range("G342").Value = 3130

# Add new row 343
$ws.Range("A343").Value = 44237
$ws.Range("A343").NumberFormat = "yyyy-mm-dd"
$ws.Range("B343").Value = 271473
$ws.Range("C343").Value = 10932
$ws.Range("D343").Value = 2487
$ws.Range("E343").Value = 5629
$ws.Range("F343").Value = 148022
$ws.Range("G343").Value = 3145
